$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $c = $Sheet.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '75.586.73'
$ws.Range("E2").Value = '  +8.61%  '

# Row 3
$ws.Range("D3").Value = '2.673.34'
$ws.Range("E3").Value = '  +9.86%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
Set-TextValue $ws "D5" '187.06'
$ws.Range("E5").Value = '  +12.44%  '

# Row 6
Set-TextValue $ws "D6" '586.69'
$ws.Range("E6").Value = '  +4.07%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("E8").Value = '  +4.89%  '

# Row 9
Set-TextValue $ws "D9" '0.194'
$ws.Range("E9").Value = '  +14.54%  '

# Row 10
$ws.Range("D10").Value = '2.673.39'
$ws.Range("E10").Value = '  +9.90%  '

# Row 11
$ws.Range("E11").Value = '  +1.89%  '

# Row 12
Set-TextValue $ws "D12" '0.358'
$ws.Range("E12").Value = '  +6.91%  '

# Row 13
$ws.Range("E13").Value = '  +1.13%  '

# Row 14
$ws.Range("D14").Value = '75.411.03'
$ws.Range("E14").Value = '  +8.59%  '

# Row 15
$ws.Range("D15").Value = '3.169.95'
$ws.Range("E15").Value = '  +10.04%  '

# Row 16
$ws.Range("E16").Value = '  +5.51%  '

# Row 17
Set-TextValue $ws "D17" '26.51'
$ws.Range("E17").Value = '  +10.58%  '

# Row 18
$ws.Range("D18").Value = '2.681.07'
$ws.Range("E18").Value = '  +11.09%  '

# Row 19
Set-TextValue $ws "D19" '9.29'
$ws.Range("E19").Value = '  +30.17%  '

# Row 20
Set-TextValue $ws "D20" '11.95'
$ws.Range("E20").Value = '  +10.56%  '

# Row 21
Set-TextValue $ws "D21" '372.27'
$ws.Range("E21").Value = '  +8.87%  '

# Row 22
Set-TextValue $ws "D22" '2.29'
$ws.Range("E22").Value = '  +16.03%  '

# Row 23
$ws.Range("E23").Value = '  +5.04%  '

# Row 24
$ws.Range("E24").Value = '  +4.20%  '

# Row 25
$ws.Range("E25").Value = '  -0.08%  '

# Row 26
Set-TextValue $ws "D26" '69.94'
$ws.Range("E26").Value = '  +6.09%  '

# Row 27
Set-TextValue $ws "D27" '4.18'
$ws.Range("E27").Value = '  +9.38%  '

# Row 28
Set-TextValue $ws "D28" '9.36'
$ws.Range("E28").Value = '  +9.86%  '

# Row 29
$ws.Range("D29").Value = '2.826.10'
$ws.Range("E29").Value = '  +10.62%  '

# Row 30
Set-TextValue $ws "D30" '1.01'
$ws.Range("E30").Value = '  +0.70%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0942'
$ws.Range("E31").Value = '  +10.64%  '

# Row 32
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws "D32" '519.92'
$ws.Range("E32").Value = '  +14.84%  '

# Row 33
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws "D33" '1.41'
$ws.Range("E33").Value = '  +13.76%  '

# Row 34
$ws.Range("E34").Value = '  +4.71%  '

# Row 35
$ws.Range("E35").Value = '  +8.66%  '

# Row 36
$ws.Range("E36").Value = '  -0.01%  '

# Row 37
Set-TextValue $ws "D37" '163.43'
$ws.Range("E37").Value = '  +3.70%  '

# Row 38
$ws.Range("E38").Value = '  +7.55%  '

# Row 39
Set-TextValue $ws "D39" '19.19'
$ws.Range("E39").Value = '  +5.29%  '

# Row 40
Set-TextValue $ws "D40" '19.37'
$ws.Range("E40").Value = '  +1.34%  '

# Row 41
$ws.Range("E41").Value = '  -0.02%  '

# Row 42
$ws.Range("E42").Value = '  +13.73%  '

# Row 43
Set-TextValue $ws "D43" '168.55'
$ws.Range("E43").Value = '  +24.76%  '

# Row 44
$ws.Range("E44").Value = '  +11.96%  '

# Row 45
Set-TextValue $ws "D45" '0.330'
$ws.Range("E45").Value = '  +9.15%  '

# Row 46
$ws.Range("E46").Value = '  +9.59%  '

# Row 47
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws "D47" '39.33'
$ws.Range("E47").Value = '  +3.89%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws "D48" '2.35'
$ws.Range("E48").Value = '  +12.26%  '

# Row 49
Set-TextValue $ws "D49" '0.0844'
$ws.Range("E49").Value = '  +16.37%  '

# Row 50
Set-TextValue $ws "D50" '3.64'
$ws.Range("E50").Value = '  +7.20%  '

# Row 51
Set-TextValue $ws "D51" '0.534'
$ws.Range("E51").Value = '  +9.30%  '
